# Fb_scenarios.xlsx — add SignUp scenario sheet + cross-browser tweak to FbLogin
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. FbLogin: insert a new step "verify login page / NA / get title / NA"
#    right before the existing "close browser" row. Rows.Insert() shifts the
#    following row down and carries its formatting, so the new row already
#    picks up the same body style used throughout the sheet.
# ---------------------------------------------------------------------------
$ws1.Rows.Item(7).Insert()
$ws1.Range("A7").Value = "verify login page"
$ws1.Range("B7").Value = "NA"
$ws1.Range("C7").Value = "get title"
$ws1.Range("D7").Value = "NA"

# ---------------------------------------------------------------------------
# 2. FbLogin: switch the browser under test from chrome to firefox.
# ---------------------------------------------------------------------------
$ws1.Range("D2").Value = "firefox"

# ---------------------------------------------------------------------------
# 3. FbLogin: the launch-url cell now shows the bare facebook URL (no
#    trailing slash) and is restyled with Excel's built-in Hyperlink look.
#    The underlying hyperlink relationship/target is left untouched.
# ---------------------------------------------------------------------------
$ws1.Range("D3").Value = "https://www.facebook.com"
$ws1.Range("D3").Style = "Hyperlink"

[void]$ws1.Range("D2").Select()

# ---------------------------------------------------------------------------
# 4. Add the new "SignUp" sheet after FbLogin, mirroring its layout.
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "SignUp"

$ws2.Columns.Item(1).ColumnWidth = 26.42578125
$ws2.Columns.Item(2).ColumnWidth = 27
$ws2.Columns.Item(3).ColumnWidth = 27.42578125
$ws2.Columns.Item(4).ColumnWidth = 27.42578125

# Header row (teststep | locator | action | value) — reuse FbLogin's bold style
$ws1.Range("A1:D1").Copy($ws2.Range("A1:D1"))

# open browser / NA / open browser / chrome — reuse FbLogin's body style
$ws1.Range("A2:D2").Copy($ws2.Range("A2:D2"))
$ws2.Range("D2").Value = "chrome"

# launch url / NA / enter url / <hyperlink> — reuse FbLogin's body style, then
# give D3 the same treatment as FbLogin!D3 (bare URL + Hyperlink style)
$ws1.Range("A3:D3").Copy($ws2.Range("A3:D3"))
$ws2.Range("D3").Value = "https://www.facebook.com"
$ws2.Range("D3").Style = "Hyperlink"

# verify sign up link / linkText=Sign Up / click / NA — typed fresh, so it
# keeps the worksheet's default (unstyled) formatting.
$ws2.Range("A4").Value = "verify sign up link"
$ws2.Range("B4").Value = "linkText=Sign Up"
$ws2.Range("C4").Value = "click"
$ws2.Range("D4").Value = "NA"

# close browser / NA / quit / NA — reuse FbLogin's body style
$ws1.Range("A7:D7").Copy($ws2.Range("A5:D5"))
$ws2.Range("A5").Value = "close browser"
$ws2.Range("B5").Value = "NA"
$ws2.Range("C5").Value = "quit"
$ws2.Range("D5").Value = "NA"

[void]$ws2.Range("D3").Select()
[void]$ws1.Activate()
